$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 82
$ws.Cells.Item(82, 8).Value = 199
$ws.Cells.Item(82, 9).Value = 199
$ws.Cells.Item(82, 11).Value = 597
$ws.Cells.Item(82, 13).Value = -191
# Row 85
$ws.Cells.Item(85, 8).Value = 199
$ws.Cells.Item(85, 9).Value = 199
$ws.Cells.Item(85, 11).Value = 597
$ws.Cells.Item(85, 13).Value = 807
# Row 112
$ws.Cells.Item(112, 8).Value = 2038.9565
$ws.Cells.Item(112, 9).Value = 1136
$ws.Cells.Item(112, 11).Value = 3408
$ws.Cells.Item(112, 13).Value = -2300
# Row 113
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 13).ClearContents()
# Row 132
$ws.Cells.Item(132, 8).Value = 1045.8108
$ws.Cells.Item(132, 9).Value = 996.90625
$ws.Cells.Item(132, 11).Value = 2990.71875
$ws.Cells.Item(132, 13).Value = -460.71875
# Row 135
$ws.Cells.Item(135, 8).Value = 1204.5454
$ws.Cells.Item(135, 9).Value = 750.1111
$ws.Cells.Item(135, 11).Value = 6750.9999
$ws.Cells.Item(135, 13).Value = -4215.9999
# Row 138
$ws.Cells.Item(138, 8).Value = 4210.0483
$ws.Cells.Item(138, 10).Value = 4525.5093
$ws.Cells.Item(138, 12).Value = 13576.5279
$ws.Cells.Item(138, 14).Value = -23856.5279

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 691.3570999999999
$ws.Cells.Item(2, 9).Value = 564.9167
$ws.Cells.Item(2, 11).Value = 564.9167
$ws.Cells.Item(2, 13).Value = -451.9167
# Row 32
$ws.Cells.Item(32, 8).Value = 14562.871
$ws.Cells.Item(32, 9).Value = 5755.3784
$ws.Cells.Item(32, 11).Value = 5755.3784
$ws.Cells.Item(32, 13).Value = -5468.3784
# Row 45
$ws.Cells.Item(45, 8).Value = 6201.75
$ws.Cells.Item(45, 9).Value = 2103
$ws.Cells.Item(45, 11).Value = 2103
$ws.Cells.Item(45, 13).Value = -1726
# Row 61
$ws.Cells.Item(61, 8).Value = 1497.5
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 13).ClearContents()
# Row 63
$ws.Cells.Item(63, 8).Value = 5423.6875
$ws.Cells.Item(63, 10).Value = 7857.143
$ws.Cells.Item(63, 12).Value = 7857.143
$ws.Cells.Item(63, 14).Value = -9229.143
# Row 66
$ws.Cells.Item(66, 8).Value = 5423.6875
$ws.Cells.Item(66, 10).Value = 7857.143
$ws.Cells.Item(66, 12).Value = 39285.715
$ws.Cells.Item(66, 14).Value = -46149.715
# Row 102
$ws.Cells.Item(102, 8).Value = 1325.8182
$ws.Cells.Item(102, 9).Value = 953.7778
$ws.Cells.Item(102, 11).Value = 953.7778
$ws.Cells.Item(102, 13).Value = 668.2222
# Row 116
$ws.Cells.Item(116, 8).Value = 691.3570999999999
$ws.Cells.Item(116, 9).Value = 564.9167
$ws.Cells.Item(116, 11).Value = 564.9167
$ws.Cells.Item(116, 13).Value = 1729.0833
# Row 136
$ws.Cells.Item(136, 8).Value = 1497.5
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 13).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 691.3570999999999
$ws.Cells.Item(3, 9).Value = 564.9167
$ws.Cells.Item(3, 11).Value = 564.9167
$ws.Cells.Item(3, 13).Value = -450.9167
# Row 86
$ws.Cells.Item(86, 8).Value = 1705.625
$ws.Cells.Item(86, 9).Value = 397.8
$ws.Cells.Item(86, 11).Value = 397.8
$ws.Cells.Item(86, 13).Value = 725.2
# Row 89
$ws.Cells.Item(89, 8).Value = 1705.625
$ws.Cells.Item(89, 9).Value = 397.8
$ws.Cells.Item(89, 11).Value = 1989
$ws.Cells.Item(89, 13).Value = 3627
# Row 105
$ws.Cells.Item(105, 8).Value = 3707.1875
$ws.Cells.Item(105, 9).Value = 2987.238
$ws.Cells.Item(105, 11).Value = 2987.238
$ws.Cells.Item(105, 13).Value = -1240.238
# Row 107
$ws.Cells.Item(107, 8).Value = 696.6667
$ws.Cells.Item(107, 10).Value = 90
$ws.Cells.Item(107, 12).Value = 90
$ws.Cells.Item(107, 14).Value = -3930
# Row 134
$ws.Cells.Item(134, 8).Value = 3335.9285
$ws.Cells.Item(134, 9).Value = 3417.2
$ws.Cells.Item(134, 11).Value = 10251.6
$ws.Cells.Item(134, 13).Value = -7716.599999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 4636.2
$ws.Cells.Item(31, 9).Value = 3069.4285
$ws.Cells.Item(31, 11).Value = 3069.4285
$ws.Cells.Item(31, 13).Value = -2774.4285
# Row 34
$ws.Cells.Item(34, 8).Value = 4636.2
$ws.Cells.Item(34, 9).Value = 3069.4285
$ws.Cells.Item(34, 11).Value = 3069.4285
$ws.Cells.Item(34, 13).Value = -2867.4285
# Row 38
$ws.Cells.Item(38, 8).Value = 38
$ws.Cells.Item(38, 9).Value = 38
$ws.Cells.Item(38, 11).Value = 38
$ws.Cells.Item(38, 13).Value = 339
# Row 46
$ws.Cells.Item(46, 8).Value = 38
$ws.Cells.Item(46, 9).Value = 38
$ws.Cells.Item(46, 11).Value = 38
$ws.Cells.Item(46, 13).Value = 173
# Row 58
$ws.Cells.Item(58, 8).Value = 3622.9412
$ws.Cells.Item(58, 9).Value = 2078.25
$ws.Cells.Item(58, 11).Value = 2078.25
$ws.Cells.Item(58, 13).Value = -1875.25
# Row 109
$ws.Cells.Item(109, 8).Value = 14322.917
$ws.Cells.Item(109, 10).Value = 14322.917
$ws.Cells.Item(109, 12).Value = 14322.917
$ws.Cells.Item(109, 14).Value = -16402.917
# Row 136
$ws.Cells.Item(136, 8).Value = 3622.9412
$ws.Cells.Item(136, 9).Value = 2078.25
$ws.Cells.Item(136, 11).Value = 6234.75
$ws.Cells.Item(136, 13).Value = -3684.75

$ws = $wb.Worksheets.Item("CUL")
# Row 59
$ws.Cells.Item(59, 8).Value = 5966
$ws.Cells.Item(59, 9).Value = 4949
$ws.Cells.Item(59, 11).Value = 14847
$ws.Cells.Item(59, 13).Value = -14307
# Row 60
$ws.Cells.Item(60, 8).Value = 454.7143
$ws.Cells.Item(60, 9).Value = 251.6
$ws.Cells.Item(60, 11).Value = 754.8
$ws.Cells.Item(60, 13).Value = -503.8
# Row 97
$ws.Cells.Item(97, 8).Value = 168
# Row 121
$ws.Cells.Item(121, 8).Value = 1243.7
$ws.Cells.Item(121, 10).Value = 2251.4
$ws.Cells.Item(121, 12).Value = 6754.200000000001
$ws.Cells.Item(121, 14).Value = -9374.200000000001
# Row 122
$ws.Cells.Item(122, 8).Value = 887
$ws.Cells.Item(122, 9).Value = 733.2857
$ws.Cells.Item(122, 11).Value = 6599.571300000001
$ws.Cells.Item(122, 13).Value = -4149.571300000001

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 3062.6667
$ws.Cells.Item(80, 9).Value = 3116.889
$ws.Cells.Item(80, 10).Value = 2900
$ws.Cells.Item(80, 11).Value = 3116.889
$ws.Cells.Item(80, 12).Value = 2900
$ws.Cells.Item(80, 13).Value = -2118.889
$ws.Cells.Item(80, 14).Value = -4896
# Row 83
$ws.Cells.Item(83, 8).Value = 3062.6667
$ws.Cells.Item(83, 9).Value = 3116.889
$ws.Cells.Item(83, 10).Value = 2900
$ws.Cells.Item(83, 11).Value = 15584.445
$ws.Cells.Item(83, 12).Value = 14500
$ws.Cells.Item(83, 13).Value = -10592.445
$ws.Cells.Item(83, 14).Value = -24484
# Row 132
$ws.Cells.Item(132, 8).Value = 3502.5264
$ws.Cells.Item(132, 9).Value = 3103.9375
$ws.Cells.Item(132, 11).Value = 9311.8125
$ws.Cells.Item(132, 13).Value = -6781.8125
# Row 136
$ws.Cells.Item(136, 8).Value = 25704.182
$ws.Cells.Item(136, 10).Value = 25704.182
$ws.Cells.Item(136, 12).Value = 77112.546
$ws.Cells.Item(136, 14).Value = -82212.546

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Cells.Item(68, 8).Value = 2994
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).ClearContents()
# Row 71
$ws.Cells.Item(71, 8).Value = 2994
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).ClearContents()
# Row 132
$ws.Cells.Item(132, 8).Value = 4787.884
$ws.Cells.Item(132, 9).Value = 4470.552
$ws.Cells.Item(132, 10).Value = 5445.2144
$ws.Cells.Item(132, 11).Value = 13411.656
$ws.Cells.Item(132, 12).Value = 16335.6432
$ws.Cells.Item(132, 13).Value = -10881.656
$ws.Cells.Item(132, 14).Value = -21395.6432
# Row 136
$ws.Cells.Item(136, 8).Value = 3395.158
$ws.Cells.Item(136, 9).Value = 3250.625
$ws.Cells.Item(136, 11).Value = 9751.875
$ws.Cells.Item(136, 13).Value = -7201.875

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 1933.1111
$ws.Cells.Item(81, 9).Value = 2199.75
$ws.Cells.Item(81, 11).Value = 4399.5
$ws.Cells.Item(81, 13).Value = -3338.5
# Row 84
$ws.Cells.Item(84, 8).Value = 1933.1111
$ws.Cells.Item(84, 9).Value = 2199.75
$ws.Cells.Item(84, 11).Value = 21997.5
$ws.Cells.Item(84, 13).Value = -16693.5
# Row 113
$ws.Cells.Item(113, 8).Value = 1393.25
$ws.Cells.Item(113, 9).Value = 1336.6364
$ws.Cells.Item(113, 10).Value = 1441.1538
$ws.Cells.Item(113, 11).Value = 4009.9092
$ws.Cells.Item(113, 12).Value = 4323.4614
$ws.Cells.Item(113, 13).Value = -1839.9092
$ws.Cells.Item(113, 14).Value = -8663.4614
# Row 132
$ws.Cells.Item(132, 8).Value = 1350.5
$ws.Cells.Item(132, 9).Value = 1297.7333
$ws.Cells.Item(132, 11).Value = 3893.199900000001
$ws.Cells.Item(132, 13).Value = -1363.199900000001
# Row 136
$ws.Cells.Item(136, 8).Value = 44874.176
$ws.Cells.Item(136, 9).Value = 873.7646999999999
$ws.Cells.Item(136, 11).Value = 2621.2941
$ws.Cells.Item(136, 13).Value = -71.29410000000007
